# Rows 5-9 on the active sheet need to be cyclically rotated:
# the last row of the block (row 9) moves up to become the new row 5,
# and rows 5-8 each shift down by one (5->6, 6->7, 7->8, 8->9).
#
# This is exactly what Excel's "Insert Cut Cells" does when you cut row 9
# and drop it on row 5. We replicate it with a sequence the COM model
# supports reliably: insert a blank row at 5 (pushing 5-9 down to 6-10),
# then cut row 10 (the original row 9, now shifted down) into that blank
# row 5, and finally remove the now-empty row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()
$ws.Range("A10:AY10").Cut($ws.Range("A5:AY5"))
$ws.Rows("10:10").Delete()
